$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing text storage (avoids Excel's
# automatic number coercion for numeric-looking strings like '20.10' or
# '0.573'), then restore the cell's original style so no spurious
# formatting diff is introduced.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '52.135.84'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '2.839.68'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E4").Value = '  -0.05%  '
Set-TextValue $ws.Range("D5") '362.51'
$ws.Range("E5").Value = '  +6.74%  '
Set-TextValue $ws.Range("D6") '113.06'
$ws.Range("E6").Value = '  -2.29%  '
Set-TextValue $ws.Range("D7") '0.573'
$ws.Range("E7").Value = '  +4.36%  '
$ws.Range("E8").Value = '  -0.03%  '
Set-TextValue $ws.Range("D9") '0.603'
$ws.Range("E9").Value = '  +4.52%  '
Set-TextValue $ws.Range("D10") '41.66'
$ws.Range("E10").Value = '  -0.60%  '
Set-TextValue $ws.Range("D11") '0.0863'
$ws.Range("E11").Value = '  +0.10%  '
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D12") '20.10'
$ws.Range("E12").Value = '  +0.05%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D13") '0.132'
$ws.Range("E13").Value = '  +1.12%  '
Set-TextValue $ws.Range("D14") '7.80'
$ws.Range("E14").Value = '  +2.54%  '
$ws.Range("D15").Value = '3.285.15'
$ws.Range("E15").Value = '  +1.94%  '
$ws.Range("D16").Value = '2.834.16'
$ws.Range("E16").Value = '  +1.77%  '
Set-TextValue $ws.Range("D17") '0.914'
$ws.Range("E17").Value = '  +3.69%  '
$ws.Range("D18").Value = '52.075.32'
$ws.Range("E18").Value = '  +0.31%  '
Set-TextValue $ws.Range("D19") '7.56'
$ws.Range("E19").Value = '  +8.70%  '
$ws.Range("E20").Value = '  -1.20%  '
Set-TextValue $ws.Range("D21") '13.58'
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("D22").Value = '0.0₂01000'
$ws.Range("E22").Value = '  +2.07%  '
Set-TextValue $ws.Range("D23") '70.37'
$ws.Range("E23").Value = '  +0.47%  '
Set-TextValue $ws.Range("D24") '269.02'
$ws.Range("E24").Value = '  -3.22%  '
$ws.Range("E25").Value = '  +4.44%  '
Set-TextValue $ws.Range("D26") '27.10'
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("E29").Value = '  +1.14%  '
$ws.Range("E30").Value = '  +28.74%  '
Set-TextValue $ws.Range("D31") '53.81'
$ws.Range("E31").Value = '  +7.22%  '
$ws.Range("E32").Value = '  -0.89%  '
Set-TextValue $ws.Range("D33") '35.57'
$ws.Range("E33").Value = '  +2.11%  '
Set-TextValue $ws.Range("D34") '5.87'
$ws.Range("E34").Value = '  +2.71%  '
Set-TextValue $ws.Range("D35") '5.46'
$ws.Range("E35").Value = '  +10.32%  '
Set-TextValue $ws.Range("D36") '0.0845'
$ws.Range("E36").Value = '  +2.50%  '
$ws.Range("E37").Value = '  -0.08%  '
$ws.Range("E38").Value = '  +1.28%  '
$ws.Range("E39").Value = '  -1.83%  '
Set-TextValue $ws.Range("D40") '18.46'
$ws.Range("E40").Value = '  -2.42%  '
$ws.Range("E41").Value = '  +1.83%  '
Set-TextValue $ws.Range("D42") '23.63'
$ws.Range("E42").Value = '  +1.99%  '
Set-TextValue $ws.Range("D43") '126.86'
$ws.Range("E43").Value = '  +1.25%  '
$ws.Range("E44").Value = '  -6.35%  '
$ws.Range("E45").Value = '  -3.19%  '
Set-TextValue $ws.Range("D46") '3.42'
$ws.Range("E46").Value = '  +3.44%  '
$ws.Range("D47").Value = '2.114.61'
$ws.Range("E47").Value = '  +1.48%  '
Set-TextValue $ws.Range("D48") '2.26'
$ws.Range("E48").Value = '  +0.97%  '
Set-TextValue $ws.Range("D49") '0.988'
$ws.Range("E49").Value = '  +11.17%  '
$ws.Range("E50").Value = '  +6.00%  '
Set-TextValue $ws.Range("D51") '62.12'
$ws.Range("E51").Value = '  +4.17%  '
